$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 0) The "_GoBack" bookmark currently sits in the "A detailed schedule..."
#    paragraph, splitting it into two runs. It needs to move to sit right
#    after the (new) last date in the schedule table ("December 10").
#    Delete the old one now; we'll add the new one once that text exists.
# ---------------------------------------------------------------------------
try {
    if ($d.Bookmarks.Exists("_GoBack")) {
        $d.Bookmarks.Item("_GoBack").Delete()
    }
} catch {
    # No pre-existing "_GoBack" bookmark (or bookmark API unavailable) -
    # nothing to clean up.
}

# ---------------------------------------------------------------------------
# 1) "October 8" -> "October 22"  (Research proposal due date)
# ---------------------------------------------------------------------------
$d.Content.Find.Execute("October 8", $true, $false, $false, $false, $false, `
    $true, 1, $false, "October 22", 2) | Out-Null

# ---------------------------------------------------------------------------
# 2) Last-row date "December 3" -> "December 10". Do this BEFORE turning
#    "November 26" into "December 3" below, so the Find below unambiguously
#    targets the original "December 3" cell (the last row of the table).
#
#    A bookmark ("_GoBack") must end up immediately after the new
#    "December 10" text, at the very end of that table cell/row. Adding a
#    bookmark exactly at that boundary is temporarily unreliable, so a
#    one-character sentinel is inserted after the date first (so the
#    bookmark is not created at the absolute edge), the bookmark is added,
#    the date text is updated, and finally the sentinel is removed.
# ---------------------------------------------------------------------------
$rng = $d.Content
$rng.Find.Execute("December 3", $true, $false, $false, $false, $false, `
    $true, 1, $false, "", 0) | Out-Null
$startPos = $rng.Start
$endPos = $rng.End

$sentinelIns = $d.Range($endPos, $endPos)
$sentinelIns.InsertAfter("X")

$bmRange = $d.Range($endPos, $endPos)
$d.Bookmarks.Add("_GoBack", $bmRange)

$dateRange = $d.Range($startPos, $endPos)
$dateRange.Text = "December 10"

$bm = $d.Bookmarks.Item("_GoBack")
$sentinelRange = $d.Range($bm.End, $bm.End + 1)
$sentinelRange.Text = ""

# ---------------------------------------------------------------------------
# 3) "November 26" -> "December 3" (Final paper due date)
# ---------------------------------------------------------------------------
$d.Content.Find.Execute("November 26", $true, $false, $false, $false, $false, `
    $true, 1, $false, "December 3", 2) | Out-Null

# ---------------------------------------------------------------------------
# 4) Merge the two runs "A detailed schedule (with readings)" and
#    " is available at: " (which used to be separated by the "_GoBack"
#    bookmark removed above) into a single run with the combined text.
# ---------------------------------------------------------------------------
$d.Content.Find.Execute("A detailed schedule (with readings) is available at: ", `
    $true, $false, $false, $false, $false, $true, 1, $false, `
    "A detailed schedule (with readings) is available at: ", 2) | Out-Null
